$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param([string]$Address, [string]$Val)
    $c = $ws.Range($Address)
    $c.NumberFormat = "@"
    $c.Value = $Val
    $c.Style = "Normal"
}

Set-TextValue "D2" "66.493.16"
Set-TextValue "E2" "  +3.68%  "
Set-TextValue "D3" "3.482.55"
Set-TextValue "E3" "  +2.26%  "
Set-TextValue "D5" "589.92"
Set-TextValue "E5" "  +2.90%  "
Set-TextValue "D6" "168.15"
Set-TextValue "E6" "  +3.27%  "
Set-TextValue "E7" "  -0.03%  "
Set-TextValue "D8" "3.480.12"
Set-TextValue "E8" "  +2.29%  "
Set-TextValue "D9" "0.587"
Set-TextValue "E9" "  +6.71%  "
Set-TextValue "D10" "7.31"
Set-TextValue "E10" "  +0.04%  "
Set-TextValue "E11" "  +5.52%  "
Set-TextValue "D12" "0.434"
Set-TextValue "E12" "  +3.04%  "
Set-TextValue "D13" "4.088.31"
Set-TextValue "E13" "  +2.41%  "
Set-TextValue "E14" "  -0.67%  "
Set-TextValue "D15" "28.04"
Set-TextValue "E15" "  +4.29%  "
Set-TextValue "D16" "66.548.56"
Set-TextValue "E16" "  +3.79%  "
Set-TextValue "D17" "0.0000176"
Set-TextValue "E17" "  +2.12%  "
Set-TextValue "D18" "3.487.50"
Set-TextValue "E18" "  +3.23%  "
Set-TextValue "E19" "  +2.79%  "
Set-TextValue "D20" "13.93"
Set-TextValue "E20" "  +3.68%  "
Set-TextValue "D21" "389.98"
Set-TextValue "E21" "  +4.07%  "
Set-TextValue "D22" "7.88"
Set-TextValue "E22" "  +1.14%  "
Set-TextValue "D23" "72.85"
Set-TextValue "E23" "  +3.46%  "
Set-TextValue "E24" "  -0.15%  "
Set-TextValue "D25" "0.531"
Set-TextValue "E25" "  +3.78%  "
Set-TextValue "D26" "0.0000120"
Set-TextValue "E26" "  +5.08%  "
Set-TextValue "D27" "10.22"
Set-TextValue "E27" "  +7.74%  "
Set-TextValue "E28" "  +1.02%  "
Set-TextValue "D30" "6.29"
Set-TextValue "E30" "  +3.75%  "
Set-TextValue "D31" "1.44"
Set-TextValue "E31" "  +3.72%  "
Set-TextValue "D32" "2.05"
Set-TextValue "E32" "  +2.67%  "
Set-TextValue "D33" "23.51"
Set-TextValue "E33" "  +3.09%  "
Set-TextValue "D34" "7.34"
Set-TextValue "E34" "  +4.25%  "
Set-TextValue "D35" "1.00"
Set-TextValue "E35" "  +0.00%  "
Set-TextValue "E36" "  +8.18%  "
Set-TextValue "D37" "162.36"
Set-TextValue "E37" "  +1.83%  "
Set-TextValue "E38" "  +3.75%  "
Set-TextValue "E39" "  +4.80%  "
Set-TextValue "D40" "6.73"
Set-TextValue "E40" "  +4.93%  "
Set-TextValue "D41" "0.0741"
Set-TextValue "E41" "  +2.08%  "
Set-TextValue "D42" "4.61"
Set-TextValue "E42" "  +5.82%  "
Set-TextValue "D43" "26.22"
Set-TextValue "E43" "  +1.71%  "
Set-TextValue "B44" "OKB"
Set-TextValue "C44" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D44" "43.05"
Set-TextValue "E44" "  +0.90%  "
Set-TextValue "B45" "InjectiveProtocol"
Set-TextValue "C45" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D45" "26.55"
Set-TextValue "E45" "  +3.11%  "
Set-TextValue "D46" "2.761.64"
Set-TextValue "E46" "  +1.37%  "
Set-TextValue "D47" "0.0310"
Set-TextValue "E47" "  +2.05%  "
Set-TextValue "D48" "2.48"
Set-TextValue "E48" "  +3.22%  "
Set-TextValue "D49" "343.82"
Set-TextValue "E49" "  +4.53%  "
Set-TextValue "E50" "  +3.74%  "
Set-TextValue "D51" "33.67"
Set-TextValue "E51" "  +11.68%  "
